# Update the standard-error figures (in parentheses) under the theta_se
# (row 4) and lambda_se (row 6) lines of the BAC test primary table with
# the finalized replicate results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# theta_se row (row 4), columns B..L = years 1983..1993
$ws.Range("B4").Value = "(0.15)"
$ws.Range("C4").Value = "(0.33)"
$ws.Range("D4").Value = "(2.79)"
$ws.Range("E4").Value = "(0.06)"
$ws.Range("F4").Value = "(0.13)"
$ws.Range("G4").Value = "(0.38)"
$ws.Range("H4").Value = "(0.38)"
$ws.Range("I4").Value = "(1.01)"
$ws.Range("J4").Value = "(0.07)"
$ws.Range("K4").Value = "(0.97)"
$ws.Range("L4").Value = "(1.61)"

# lambda_se row (row 6), columns B..L = years 1983..1993
$ws.Range("B6").Value = "(0.08)"
$ws.Range("C6").Value = "(0.05)"
$ws.Range("D6").Value = "(1.81)"
$ws.Range("E6").Value = "(0.51)"
$ws.Range("F6").Value = "(0.5)"
$ws.Range("G6").Value = "(0.77)"
$ws.Range("H6").Value = "(0.38)"
$ws.Range("I6").Value = "(0.85)"
$ws.Range("J6").Value = "(1.3)"
$ws.Range("K6").Value = "(0.52)"
$ws.Range("L6").Value = "(1.74)"
